$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The source data table lists two rows of numeric-looking values stored
# as text for "Alex Carey" (runs/balls/fours/sixes). The activity has
# been updated so the figures previously on row 3 now belong to row 2,
# and vice versa (row 4 is unaffected).
#
# Row 2 (runs=4, balls=7, sixes=0) -> (runs=14, balls=13, sixes=1)
# Row 3 (runs=14, balls=13, sixes=1) -> (runs=4, balls=7, sixes=0)
#
# A leading apostrophe keeps these numeric-looking values stored as text,
# matching the original "number stored as text" cell formatting.
$ws.Range("C2").Value = "'14"
$ws.Range("D2").Value = "'13"
$ws.Range("F2").Value = "'1"

$ws.Range("C3").Value = "'4"
$ws.Range("D3").Value = "'7"
$ws.Range("F3").Value = "'0"
